$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits after the
#    first paragraph ("步骤：安装油猴插件，添加并启用脚本，打开阅读汇总帖。").
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) Update the wording of the "如果f12中看见..." sentence.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("被拦截的话点始终允许即可。", $false, $false, $false, $false, $false, $true, 1, $false, "被拦截的话点始终允许，再刷新一下汇总贴或者重启浏览器即可。", 2)

# $r now covers the freshly-inserted replacement text; use it to find the
# paragraph that holds the sentence so we can locate the end of it (the
# position right before the paragraph mark).
$para = $r.Paragraphs(1)
$endPos = $para.Range.End - 1

# ------------------------------------------------------------------
# 3) Re-insert a "_GoBack" bookmark, collapsed, right at the end of that
#    paragraph's text (after the final "。", before the paragraph mark).
#
#    Bookmarks.Add() in this host always snaps to the *start* of the
#    supplied Range, and a zero-length Range at a paragraph boundary gets
#    mis-resolved -- so a temporary one-character placeholder is used to
#    give Add() a safe, unambiguous, non-boundary Range to anchor on, and
#    the placeholder is then deleted, leaving a clean collapsed bookmark.
# ------------------------------------------------------------------
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos + 1)
$bmRange.Bookmarks.Add("_GoBack")

$d.Range($endPos, $endPos + 1).Delete()
